$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.377
$ws.Range("B8").Value = 5.912
$ws.Range("B10").Value = 6.148
$ws.Range("C11").Value = -12.266
$ws.Range("B12").Value = 5.258
$ws.Range("C12").Value = -10.95
$ws.Range("C15").Value = -13.364
$ws.Range("C17").Value = -13.267
$ws.Range("B18").Value = 5.529999999999999
$ws.Range("B25").Value = 5.482
$ws.Range("C26").Value = -13.15
$ws.Range("C27").Value = -13.383
$ws.Range("C28").Value = -12.77
$ws.Range("C32").Value = -12.954
$ws.Range("B37").Value = 8.190999999999999
$ws.Range("C37").Value = -12.055
$ws.Range("C41").Value = -12.375
$ws.Range("C47").Value = -12.605
$ws.Range("C51").Value = -11.371
$ws.Range("B55").Value = 4.864999999999999
$ws.Range("C65").Value = -12.166
$ws.Range("B68").Value = 4.755
$ws.Range("C73").Value = -12.344
$ws.Range("B77").Value = 6.055999999999999
$ws.Range("B78").Value = 7.45
$ws.Range("B79").Value = 5.566
$ws.Range("B80").Value = 7.581999999999999
$ws.Range("B81").Value = 5.983
$ws.Range("B82").Value = 5.621
$ws.Range("B84").Value = 5.823
$ws.Range("C84").Value = -13.055
$ws.Range("C85").Value = -12.527
$ws.Range("C89").Value = -13.483
$ws.Range("C93").Value = -10.438
$ws.Range("C95").Value = -12.451
$ws.Range("C98").Value = -13.276
$ws.Range("C99").Value = -11.696
$ws.Range("B101").Value = 6.214
$ws.Range("C101").Value = -12.715
$ws.Range("B102").Value = 5.893
$ws.Range("C102").Value = -12.611
